$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$cellData = @(
    @("75 x 25", "  2    5", "  ----", "7|    |", "5|    |"),
    @("69 x 35", "  3    5", "  ----", "6|    |", "9|    |"),
    @("17 x 47", "  4    7", "  ----", "1|    |", "7|    |"),
    @("83 x 67", "  6    7", "  ----", "8|    |", "3|    |"),
    @("73 x 48", "  4    8", "  ----", "7|    |", "3|    |"),
    @("37 x 62", "  6    2", "  ----", "3|    |", "7|    |"),
    @("17 x 72", "  7    2", "  ----", "1|    |", "7|    |"),
    @("93 x 92", "  9    2", "  ----", "9|    |", "3|    |"),
    @("60 x 64", "  6    4", "  ----", "6|    |", "0|    |"),
    @("59 x 82", "  8    2", "  ----", "5|    |", "9|    |"),
    @("38 x 12", "  1    2", "  ----", "3|    |", "8|    |"),
    @("63 x 95", "  9    5", "  ----", "6|    |", "3|    |"),
    @("60 x 68", "  6    8", "  ----", "6|    |", "0|    |"),
    @("19 x 40", "  4    0", "  ----", "1|    |", "9|    |"),
    @("47 x 28", "  2    8", "  ----", "4|    |", "7|    |"),
)

$nCols = 3
for ($i = 0; $i -lt $cellData.Count; $i++) {
    $rowIdx = [int]([math]::Floor($i / $nCols)) + 1
    $colIdx = ($i % $nCols) + 1
    $cell = $t.Rows.Item($rowIdx).Cells.Item($colIdx)
    $parts = $cellData[$i]
    $newText = [string]::Join([string][char]11, $parts)
    $cell.Range.Text = $newText
}

Write-Output "Updated $($cellData.Count) cells"
